# attendance new route added
# Applies updated attendance / PF-ESI / wage-calculation figures to the
# Mar-2024 attendance workbook (22SNCJO-373) following a revised route
# (reduced mandays per employee).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ATTENDANCE")
$ws2 = $wb.Worksheets.Item("PF_ESI")
$ws3 = $wb.Worksheets.Item("Wage_Calculation")

# ---------------------------------------------------------------------
# Sheet 1: ATTENDANCE  (daily attendance marks for March 2024)
# ---------------------------------------------------------------------

# Row 9 - NITYA SUNDAR MUDULI (SKILLED)
$ws1.Range("D9").Value  = 0
$ws1.Range("E9").Value  = 0.5
$ws1.Range("G9").Value  = 0.5
$ws1.Range("H9").Value  = 0.5
$ws1.Range("I9").Value  = 0.5
$ws1.Range("J9").Value  = 0.5
$ws1.Range("K9").Value  = 0
$ws1.Range("L9").Value  = 0.5
$ws1.Range("N9").Value  = 0.5
$ws1.Range("O9").Value  = 0.5
$ws1.Range("P9").Value  = 0.5
$ws1.Range("Q9").Value  = 0.5
$ws1.Range("R9").Value  = 0.5
$ws1.Range("S9").Value  = 0.5
$ws1.Range("U9").Value  = 0
$ws1.Range("V9").Value  = 0
$ws1.Range("W9").Value  = 0
$ws1.Range("X9").Value  = 0
$ws1.Range("Y9").Value  = 0
$ws1.Range("Z9").Value  = 0
$ws1.Range("AB9").Value = 0.5
$ws1.Range("AC9").Value = 0
$ws1.Range("AD9").Value = 0
$ws1.Range("AE9").Value = 0.5
$ws1.Range("AF9").Value = 0.5
$ws1.Range("AG9").Value = 0.5
$ws1.Range("AI9").Value = 8

# Row 10 - JAGANNATH SAHU (SEMI-SKILLED)
$ws1.Range("D10").Value  = 0
$ws1.Range("E10").Value  = 0
$ws1.Range("H10").Value  = 0
$ws1.Range("O10").Value  = 0
$ws1.Range("P10").Value  = 0
$ws1.Range("Q10").Value  = 0
$ws1.Range("R10").Value  = 0
$ws1.Range("S10").Value  = 0
$ws1.Range("T10").Value  = 1
$ws1.Range("U10").Value  = 0.5
$ws1.Range("V10").Value  = 0
$ws1.Range("W10").Value  = 0
$ws1.Range("X10").Value  = 0
$ws1.Range("Y10").Value  = 0
$ws1.Range("Z10").Value  = 0
$ws1.Range("AC10").Value = 0
$ws1.Range("AD10").Value = 0
$ws1.Range("AE10").Value = 1.5
$ws1.Range("AF10").Value = 0
$ws1.Range("AG10").Value = 0
$ws1.Range("AH10").Value = 1
$ws1.Range("AI10").Value = 11

# Row 11 - DANA MAJHI (UNSKILLED)
$ws1.Range("P11").Value  = 0
$ws1.Range("T11").Value  = 0.5
$ws1.Range("W11").Value  = 0
$ws1.Range("Z11").Value  = 0
$ws1.Range("AC11").Value = 0
$ws1.Range("AF11").Value = 0.5
$ws1.Range("AG11").Value = 0
$ws1.Range("AI11").Value = 21

# Row 13 - TOTAL
$ws1.Range("AI13").Value = 40

# ---------------------------------------------------------------------
# Sheet 2: PF_ESI  (PF and ESI calculation details)
# ---------------------------------------------------------------------

# Row 6 / Row 7 - NITYA SUNDAR MUDULI / SKILLED TOTAL
$ws2.Range("C6").Value = 8
$ws2.Range("F6").Value = 5672
$ws2.Range("H6").Value = 5672
$ws2.Range("I6").Value = 680.64
$ws2.Range("J6").Value = 42.54
$ws2.Range("L6").Value = 723.1799999999999
$ws2.Range("M6").Value = 709
$ws2.Range("N6").Value = 28.36
$ws2.Range("O6").Value = 737.36
$ws2.Range("P6").Value = 184.34
$ws2.Range("Q6").Value = 921.7
$ws2.Range("R6").Value = 5657.82

$ws2.Range("C7").Value = 8
$ws2.Range("F7").Value = 5672
$ws2.Range("H7").Value = 5672
$ws2.Range("I7").Value = 680.64
$ws2.Range("J7").Value = 42.54
$ws2.Range("L7").Value = 723.1799999999999
$ws2.Range("M7").Value = 709
$ws2.Range("N7").Value = 28.36
$ws2.Range("O7").Value = 737.36
$ws2.Range("P7").Value = 184.34
$ws2.Range("Q7").Value = 921.7
$ws2.Range("R7").Value = 5657.82

# Row 8 / Row 9 - JAGANNATH SAHU / SEMI-SKILLED TOTAL
$ws2.Range("C8").Value = 11
$ws2.Range("F8").Value = 6479
$ws2.Range("H8").Value = 6479
$ws2.Range("I8").Value = 777.48
$ws2.Range("J8").Value = 48.59
$ws2.Range("L8").Value = 826.0700000000001
$ws2.Range("M8").Value = 809.88
$ws2.Range("N8").Value = 32.4
$ws2.Range("O8").Value = 842.28
$ws2.Range("P8").Value = 210.57
$ws2.Range("Q8").Value = 1052.85
$ws2.Range("R8").Value = 6241.93

$ws2.Range("C9").Value = 11
$ws2.Range("F9").Value = 6479
$ws2.Range("H9").Value = 6479
$ws2.Range("I9").Value = 777.48
$ws2.Range("J9").Value = 48.59
$ws2.Range("L9").Value = 826.0700000000001
$ws2.Range("M9").Value = 809.88
$ws2.Range("N9").Value = 32.4
$ws2.Range("O9").Value = 842.28
$ws2.Range("P9").Value = 210.57
$ws2.Range("Q9").Value = 1052.85
$ws2.Range("R9").Value = 6241.93

# Row 10 / Row 11 - DANA MAJHI / UNSKILLED TOTAL
$ws2.Range("C10").Value = 21
$ws2.Range("F10").Value = 10584
$ws2.Range("H10").Value = 10584
$ws2.Range("I10").Value = 1270.08
$ws2.Range("J10").Value = 79.38
$ws2.Range("L10").Value = 1349.46
$ws2.Range("M10").Value = 1323
$ws2.Range("N10").Value = 52.92
$ws2.Range("O10").Value = 1375.92
$ws2.Range("P10").Value = 343.98
$ws2.Range("Q10").Value = 1719.9
$ws2.Range("R10").Value = 9738.540000000001

$ws2.Range("C11").Value = 21
$ws2.Range("F11").Value = 10584
$ws2.Range("H11").Value = 10584
$ws2.Range("I11").Value = 1270.08
$ws2.Range("J11").Value = 79.38
$ws2.Range("L11").Value = 1349.46
$ws2.Range("M11").Value = 1323
$ws2.Range("N11").Value = 52.92
$ws2.Range("O11").Value = 1375.92
$ws2.Range("P11").Value = 343.98
$ws2.Range("Q11").Value = 1719.9
$ws2.Range("R11").Value = 9738.540000000001

# Row 12 - TOTAL
$ws2.Range("C12").Value = 40
$ws2.Range("F12").Value = 22735
$ws2.Range("H12").Value = 22735
$ws2.Range("I12").Value = 2728.2
$ws2.Range("J12").Value = 170.51
$ws2.Range("L12").Value = 2898.71
$ws2.Range("M12").Value = 2841.88
$ws2.Range("N12").Value = 113.68
$ws2.Range("O12").Value = 2955.56
$ws2.Range("P12").Value = 738.89
$ws2.Range("Q12").Value = 3694.45
$ws2.Range("R12").Value = 21638.29

# ---------------------------------------------------------------------
# Sheet 3: Wage_Calculation  (Annexure-2 wage summary)
# ---------------------------------------------------------------------

# Unskilled / Semiskilled / Skilled mandays & gross salary
$ws3.Range("D8").Value  = 21
$ws3.Range("E8").Value  = 21
$ws3.Range("F8").Value  = 10584

$ws3.Range("D9").Value  = 11
$ws3.Range("E9").Value  = 11
$ws3.Range("F9").Value  = 6479

$ws3.Range("D10").Value = 8
$ws3.Range("E10").Value = 8
$ws3.Range("F10").Value = 5672

# PF / ESI / EDLI
$ws3.Range("F16").Value = 2841.88
$ws3.Range("F17").Value = 738.89
$ws3.Range("F18").Value = 113.68

# Service charges, GST, Totals
$ws3.Range("F20").Value = 886.41
$ws3.Range("F21").Value = 5241.2148
$ws3.Range("F22").Value = 34359.0748
$ws3.Range("F28").Value = 34359.0748

Write-Host "Attendance, PF/ESI and wage calculation sheets updated for revised route."
